$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2/B3 values and clear C2/C3 (which held "Ovaliderad")
$ws.Range("B2").Value = 57881
$ws.Range("C2").ClearContents()

$ws.Range("B3").Value = 57881
$ws.Range("C3").ClearContents()
